# Auto update Excel log: append 5 new mmWave sensor readings (rows 41-45)
# to the "mmWave" sheet, extending the logged range from A1:F40 to A1:F45.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$rows = @(
    @("2026-01-31", "22:00:13", "22:00", "Living Room", "NO_MOTION_DETECTED", "Inactive"),
    @("2026-01-31", "22:00:23", "22:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-31", "22:00:34", "22:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-31", "22:00:44", "22:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-31", "22:00:55", "22:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$startRow = 41
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le 6; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        # Force text entry so date/time-shaped strings (e.g. "2026-01-31")
        # are kept as literal text instead of being auto-parsed into date
        # serials, then strip the temporary "@" number format again so the
        # new cells end up with the same default (unstyled) formatting as
        # the rest of the log.
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$c - 1]
        $cell.ClearFormats()
    }
}
